# Repull data, push all data, mean calculation
# Update the dSF column (F) values to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -5
$ws.Range("F6").Value = -8
$ws.Range("F8").Value = -3
$ws.Range("F11").Value = -3
$ws.Range("F17").Value = 1
$ws.Range("F20").Value = -2
$ws.Range("F21").Value = 7
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = -2
$ws.Range("F31").Value = -5
$ws.Range("F33").Value = -4
$ws.Range("F34").Value = 16
